# 5.31 Add SoundEffect & Optimise Save Function
#
# Adds two new trailing columns to the Sheet1 header row:
#   S1 = "SoundEffect"
#   T1 = "SEAction"
# and moves the viewport/selection so the new columns are visible
# (topLeftCell D1 -> J1, selection R10 -> U4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (plain shared-string values, same formatting as the
# rest of row 1 - no extra style applied).
$ws.Range("S1").Value = "SoundEffect"
$ws.Range("T1").Value = "SEAction"

# Scroll the view so column J is the left-most visible column, then move
# the active selection to U4 (mirrors the author re-positioning the
# viewport after adding the new columns).
$excel.ActiveWindow.ScrollColumn = 10
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("U4").Select()
